$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header text:
#  G1: "Vehicle Standby" -> "Vehicle Standby Arrival"
#  J1: "Vehicle Standby" -> "Vehicle Standby Departure" (new distinct string)
$ws.Range("G1").Value = "Vehicle Standby Arrival"
$ws.Range("J1").Value = "Vehicle Standby Departure"

# Move the active selection / view to reflect the edited area
$ws.Range("J2").Select()
$window = $excel.ActiveWindow
$window.ScrollColumn = 4
